# Atualização automática: 2025-08-26 09:00:46
# Updates detection records (image filenames, bounding-box coords and
# confidence scores) on Sheet1, plus swaps two Fly_ID rows (21 <-> 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns I (First_Coords) and J (First_Confidence) store plain text that
# looks numeric ("643,531,686,575" / "0.76"); force Text format so Excel
# doesn't silently convert them to numbers and drop formatting (e.g. the
# trailing zero in "0.70").

# Row 16
$ws.Range("D16").Value = "image_20250807111026_ppp0.jpg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "641,529,688,576"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.75"

# Row 17
$ws.Range("D17").Value = "image_20250807111026_ppp0.jpg"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "793,481,831,526"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "0.70"

# Row 18
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "1182,405,1231,455"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "0.76"

# Row 21 <-> Row 22: Fly_ID values swap places, and each row gets its own
# updated image/coords/confidence values.
$ws.Range("A21").Value = "66efa766-1456-4beb-b92a-0615a2fc41bb"
$ws.Range("D21").Value = "image_20250824214658_ppp0.jpg"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "1272,293,1315,331"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "0.69"

$ws.Range("A22").Value = "a2ea21b8-7dce-4e6a-be35-4edaddca5896"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "999,782,1040,825"
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "0.58"
